$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 161; $r++) {
    $ws.Cells.Item($r, 4).Value = 2000 + ($r - 1) * 50
}

$ws.Range("D161").Select()
